# Daily attendance processing - 2026-01-25 09:34:59
# Swap the order of the "Recorded By" contributors in column G so that the
# human recorder's email is listed before "System" (was "System, <email>",
# now "<email>, System") for every session row that still lists them in the
# old order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldOrder = "System, dnasr281@gmail.com"
$newOrder = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldOrder) {
        $cell.Value = $newOrder
    }
}
